# Update "想去人数" (interested-count) values in column F across multiple
# sheets of the 上海-漫展信息 workbook, per the generated-output refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 6597
$ws.Range("F3").Value  = 769
$ws.Range("F5").Value  = 116
$ws.Range("F6").Value  = 635
$ws.Range("F9").Value  = 809
$ws.Range("F10").Value = 1270
$ws.Range("F18").Value = 702
$ws.Range("F19").Value = 430
$ws.Range("F22").Value = 1100
$ws.Range("F23").Value = 213
$ws.Range("F24").Value = 2294
$ws.Range("F25").Value = 261
$ws.Range("F28").Value = 60
$ws.Range("F29").Value = 3728

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value  = 35
$ws.Range("F6").Value  = 733
$ws.Range("F17").Value = 389
$ws.Range("F24").Value = 215
$ws.Range("F26").Value = 100

# Sheet "本地生活"
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 1607
$ws.Range("F8").Value = 914

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 1607
$ws.Range("F7").Value  = 914
$ws.Range("F8").Value  = 6597
$ws.Range("F9").Value  = 35
$ws.Range("F10").Value = 769
$ws.Range("F11").Value = 733
$ws.Range("F12").Value = 116
$ws.Range("F13").Value = 635
$ws.Range("F16").Value = 809
$ws.Range("F23").Value = 1270
$ws.Range("F29").Value = 389
$ws.Range("F35").Value = 702
$ws.Range("F36").Value = 430
$ws.Range("F39").Value = 215
$ws.Range("F40").Value = 1100
$ws.Range("F41").Value = 213
$ws.Range("F42").Value = 2295
$ws.Range("F47").Value = 3728
